$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 710.4
$ws.Range("I2").Value = 441.54544
$ws.Range("J2").Value = 1449.75
$ws.Range("K2").Value = 441.54544
$ws.Range("L2").Value = 1449.75
$ws.Range("M2").Value = -328.54544
$ws.Range("N2").Value = -1675.75

$ws.Range("H58").Value = 438.27274
$ws.Range("I58").Value = 312.83334
$ws.Range("J58").Value = 588.8
$ws.Range("K58").Value = 938.5000200000001
$ws.Range("L58").Value = 1766.4
$ws.Range("M58").Value = -788.5000200000001
$ws.Range("N58").Value = -2066.4

$ws.Range("H86").Value = 4017.4
$ws.Range("I86").Value = 2074.6667
$ws.Range("K86").Value = 2074.6667
$ws.Range("M86").Value = -951.6667000000002

$ws.Range("H89").Value = 4017.4
$ws.Range("I89").Value = 2074.6667
$ws.Range("K89").Value = 10373.3335
$ws.Range("M89").Value = -4757.333500000001

$ws.Range("H92").Value = 4209
$ws.Range("I92").Value = 3807
$ws.Range("K92").Value = 3807
$ws.Range("M92").Value = -2559

$ws.Range("H135").Value = 3415.7778
$ws.Range("I135").Value = 1096.1428
$ws.Range("K135").Value = 9865.2852
$ws.Range("M135").Value = -7330.2852

$ws.Range("H138").Value = 2975.5
$ws.Range("I138").Value = 1894.5
$ws.Range("J138").Value = 3229.853
$ws.Range("K138").Value = 5683.5
$ws.Range("L138").Value = 9689.559000000001
$ws.Range("M138").Value = -543.5
$ws.Range("N138").Value = -19969.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 10186.75
$ws.Range("I13").Value = 17774.75
$ws.Range("J13").Value = 2598.75
$ws.Range("K13").Value = 17774.75
$ws.Range("L13").Value = 2598.75
$ws.Range("M13").Value = -17630.75
$ws.Range("N13").Value = -2886.75

$ws.Range("H32").Value = 1156.2354
$ws.Range("I32").Value = 1215.3334
$ws.Range("K32").Value = 1215.3334
$ws.Range("M32").Value = -928.3334

$ws.Range("H61").Value = 3869.5715
$ws.Range("I61").Value = 2796.5715
$ws.Range("J61").Value = 4942.5713
$ws.Range("K61").Value = 2796.5715
$ws.Range("L61").Value = 4942.5713
$ws.Range("M61").Value = -2584.5715
$ws.Range("N61").Value = -5366.5713

$ws.Range("H74").Value = 2250.2368
$ws.Range("I74").Value = 1219.6
$ws.Range("J74").Value = 6115.125
$ws.Range("K74").Value = 1219.6
$ws.Range("L74").Value = 6115.125
$ws.Range("M74").Value = -345.5999999999999
$ws.Range("N74").Value = -7863.125

$ws.Range("H77").Value = 2250.2368
$ws.Range("I77").Value = 1219.6
$ws.Range("J77").Value = 6115.125
$ws.Range("K77").Value = 6098
$ws.Range("L77").Value = 30575.625
$ws.Range("M77").Value = -1730
$ws.Range("N77").Value = -39311.625

$ws.Range("H122").Value = 1645.1428
$ws.Range("I122").Value = 1439.5264
$ws.Range("K122").Value = 4318.5792
$ws.Range("M122").Value = -1868.5792

$ws.Range("H136").Value = 3869.5715
$ws.Range("I136").Value = 2796.5715
$ws.Range("J136").Value = 4942.5713
$ws.Range("K136").Value = 8389.7145
$ws.Range("L136").Value = 14827.7139
$ws.Range("M136").Value = -5839.7145
$ws.Range("N136").Value = -19927.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2015.1538
$ws.Range("I20").Value = 1863.8572
$ws.Range("J20").Value = 2191.6667
$ws.Range("K20").Value = 1863.8572
$ws.Range("L20").Value = 2191.6667
$ws.Range("M20").Value = -1616.8572
$ws.Range("N20").Value = -2685.6667

$ws.Range("H94").Value = 1123.8422
$ws.Range("I94").Value = 1018.8333
$ws.Range("K94").Value = 1018.8333
$ws.Range("M94").Value = -567.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 30993.334
$ws.Range("J124").Value = 30993.334
$ws.Range("L124").Value = 30993.334
$ws.Range("N124").Value = -35903.334

$ws.Range("H134").Value = 2920.125
$ws.Range("I134").Value = 2766.2856
$ws.Range("K134").Value = 8298.856800000001
$ws.Range("M134").Value = -5763.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2000001
$ws.Range("I14").Value = 2000001
$ws.Range("K14").Value = 6000003
$ws.Range("M14").Value = -5999830

$ws.Range("H86").Value = 249.66667
$ws.Range("I86").Value = 249.66667
$ws.Range("K86").Value = 749.00001
$ws.Range("M86").Value = 436.99999

$ws.Range("H89").Value = 249.66667
$ws.Range("I89").Value = 249.66667
$ws.Range("K89").Value = 2247.00003
$ws.Range("M89").Value = 3680.99997

$ws.Range("H103").Value = 2560.182
$ws.Range("J103").Value = 3152.3333
$ws.Range("L103").Value = 9456.999899999999
$ws.Range("N103").Value = -11214.9999

$ws.Range("H131").Value = 1878.9166
$ws.Range("J131").Value = 2074.7
$ws.Range("L131").Value = 6224.099999999999
$ws.Range("N131").Value = -16304.1

$ws.Range("H132").Value = 1745.9333
$ws.Range("I132").Value = 1821.8182
$ws.Range("K132").Value = 16396.3638
$ws.Range("M132").Value = -13866.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 6278.8184
$ws.Range("J9").Value = 15726.25
$ws.Range("L9").Value = 15726.25
$ws.Range("N9").Value = -16066.25

$ws.Range("H97").Value = 1426.68
$ws.Range("I97").Value = 1401.8
$ws.Range("K97").Value = 1401.8
$ws.Range("M97").Value = -905.8

$ws.Range("H122").Value = 49922.668
$ws.Range("I122").Value = 55156.85
$ws.Range("K122").Value = 165470.55
$ws.Range("M122").Value = -163020.55

$ws.Range("H132").Value = 13527467
$ws.Range("I132").Value = 19244704
$ws.Range("K132").Value = 57734112
$ws.Range("M132").Value = -57731582

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6156.5557
$ws.Range("I16").Value = 8968.5
$ws.Range("J16").Value = 532.6667
$ws.Range("K16").Value = 8968.5
$ws.Range("L16").Value = 532.6667
$ws.Range("M16").Value = -8798.5
$ws.Range("N16").Value = -872.6667

$ws.Range("H22").Value = 757.5
$ws.Range("I22").Value = 697.125
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 697.125
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -402.125
$ws.Range("N22").Value = -1589

$ws.Range("H27").Value = 757.5
$ws.Range("I27").Value = 697.125
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 697.125
$ws.Range("L27").Value = 999
$ws.Range("M27").Value = -590.125
$ws.Range("N27").Value = -1213

$ws.Range("H136").Value = 5223.4194
$ws.Range("I136").Value = 2697.8235
$ws.Range("K136").Value = 8093.470499999999
$ws.Range("M136").Value = -5543.470499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4300.364
$ws.Range("I14").Value = 4144.8887
$ws.Range("K14").Value = 4144.8887
$ws.Range("M14").Value = -3976.8887

$ws.Range("H99").Value = 69420
$ws.Range("I99").Value = 69420
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 69420
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("M99").Value = -66425

$ws.Range("H122").Value = 6461.222
$ws.Range("I122").Value = 6461.222
$ws.Range("K122").Value = 19383.666
$ws.Range("M122").Value = -16933.666

Write-Output "Applied all cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
